# Applies the "Updated cryptos list" data refresh described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.931.20"
$ws.Range("E2").Value = "  +4.41%  "
$ws.Range("D3").Value = "'2.235.73"
$ws.Range("E3").Value = "  +4.61%  "
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").Value = "'251.70"
$ws.Range("E5").Value = "  +7.03%  "
$ws.Range("D6").Value = "'0.614"
$ws.Range("E6").Value = "  +2.96%  "
$ws.Range("D7").Value = "'74.61"
$ws.Range("E7").Value = "  +8.13%  "
$ws.Range("E8").Value = "  -0.23%  "
$ws.Range("D9").Value = "'0.599"
$ws.Range("E9").Value = "  +6.33%  "
$ws.Range("D10").Value = "'41.59"
$ws.Range("E10").Value = "  +8.55%  "
$ws.Range("D11").Value = "'0.0924"
$ws.Range("E11").Value = "  +4.19%  "
$ws.Range("D12").Value = "'6.89"
$ws.Range("E12").Value = "  +5.44%  "
$ws.Range("D13").Value = "'0.101"
$ws.Range("E13").Value = "  +2.21%  "
$ws.Range("D14").Value = "'2.572.70"
$ws.Range("E14").Value = "  +4.49%  "
$ws.Range("D15").Value = "'14.52"
$ws.Range("E15").Value = "  +2.52%  "
$ws.Range("D16").Value = "'2.234.85"
$ws.Range("E16").Value = "  +4.89%  "
$ws.Range("D17").Value = "'0.788"
$ws.Range("E17").Value = "  +2.14%  "
$ws.Range("D18").Value = "'42.865.21"
$ws.Range("E18").Value = "  +4.51%  "
$ws.Range("E19").Value = "  +5.05%  "
$ws.Range("D20").Value = "'71.26"
$ws.Range("E20").Value = "  +3.88%  "
$ws.Range("D21").Value = "'5.95"
$ws.Range("E21").Value = "  +4.69%  "
$ws.Range("D22").Value = "'229.71"
$ws.Range("E22").Value = "  +2.53%  "
$ws.Range("B23").Value = "ImmutableX"
$ws.Range("C23").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D23").Value = "'2.18"
$ws.Range("E23").Value = "  +13.72%  "
$ws.Range("B24").Value = "InternetComputer(DFINITY)"
$ws.Range("C24").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D24").Value = "'9.48"
$ws.Range("E24").Value = "  +0.37%  "
$ws.Range("E25").Value = "  -0.07%  "
$ws.Range("D26").Value = "'10.68"
$ws.Range("E26").Value = "  +2.24%  "
$ws.Range("D27").Value = "'3.44"
$ws.Range("E27").Value = "  +6.13%  "
$ws.Range("D28").Value = "'38.68"
$ws.Range("E28").Value = "  +24.98%  "
$ws.Range("D29").Value = "'2.21"
$ws.Range("E29").Value = "  +3.58%  "
$ws.Range("D30").Value = "'2.14"
$ws.Range("E30").Value = "  -0.37%  "
$ws.Range("D31").Value = "'171.47"
$ws.Range("E31").Value = "  +0.85%  "
$ws.Range("D32").Value = "'20.19"
$ws.Range("E32").Value = "  +3.65%  "
$ws.Range("D33").Value = "'0.0795"
$ws.Range("E33").Value = "  +6.83%  "
$ws.Range("D34").Value = "'5.21"
$ws.Range("E34").Value = "  +3.67%  "
$ws.Range("D35").Value = "'0.121"
$ws.Range("E35").Value = "  +1.94%  "
$ws.Range("E36").Value = "  +7.96%  "
$ws.Range("D37").Value = "'4.44"
$ws.Range("E37").Value = "  +7.95%  "
$ws.Range("D38").Value = "'0.0328"
$ws.Range("E38").Value = "  +15.90%  "
$ws.Range("D39").Value = "'12.63"
$ws.Range("E39").Value = "  +9.06%  "
$ws.Range("D40").Value = "'2.10"
$ws.Range("E40").Value = "  +3.70%  "
$ws.Range("D41").Value = "'0.205"
$ws.Range("E41").Value = "  +11.40%  "
$ws.Range("D42").Value = "'5.38"
$ws.Range("E42").Value = "  +3.51%  "
$ws.Range("D43").Value = "'59.61"
$ws.Range("E43").Value = "  +4.50%  "
$ws.Range("D44").Value = "'8.70"
$ws.Range("E44").Value = "  +7.01%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "'102.95"
$ws.Range("E45").Value = "  +6.79%  "
$ws.Range("B46").Value = "WOONetwork"
$ws.Range("C46").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D46").Value = "'0.480"
$ws.Range("E46").Value = "  +32.27%  "
$ws.Range("D47").Value = "'0.0988"
$ws.Range("E47").Value = "  +3.64%  "
$ws.Range("D48").Value = "'2.45"
$ws.Range("E48").Value = "  +15.73%  "
$ws.Range("E49").Value = "  +4.41%  "
$ws.Range("D50").Value = "'1.14"
$ws.Range("E50").Value = "  +4.26%  "
$ws.Range("E51").Value = "  +2.84%  "
